# Update crypto price/volume figures per the latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.621.57'
$ws.Range("E2").Value = '  +1.96%  '
# Row 3
$ws.Range("D3").Value = '1.578.55'
$ws.Range("E3").Value = '  +0.15%  '
# Row 4
$ws.Range("E4").Value = '  +0.18%  '
# Row 5
$ws.Range("D5").Value = '''212.40'
$ws.Range("E5").Value = '  -0.21%  '
# Row 6
$ws.Range("E6").Value = '  +0.02%  '
# Row 7
$ws.Range("E7").Value = '  +0.16%  '
# Row 8
$ws.Range("D8").Value = '''46.83'
$ws.Range("E8").Value = '  +7.56%  '
# Row 9
$ws.Range("D9").Value = '''24.26'
$ws.Range("E9").Value = '  +4.50%  '
# Row 10
$ws.Range("E10").Value = '  -1.28%  '
# Row 11
$ws.Range("E11").Value = '  -1.07%  '
# Row 12
$ws.Range("E12").Value = '  -0.02%  '
# Row 13
$ws.Range("D13").Value = '1.805.38'
$ws.Range("E13").Value = '  +0.30%  '
# Row 14
$ws.Range("D14").Value = '1.572.72'
$ws.Range("E14").Value = '  -0.08%  '
# Row 15
$ws.Range("E15").Value = '  +0.13%  '
# Row 16
$ws.Range("D16").Value = '''3.70'
$ws.Range("E16").Value = '  -1.39%  '
# Row 17
$ws.Range("D17").Value = '28.592.81'
$ws.Range("E17").Value = '  +1.93%  '
# Row 18
$ws.Range("D18").Value = '''62.32'
$ws.Range("E18").Value = '  -1.97%  '
# Row 19
$ws.Range("D19").Value = '''228.89'
$ws.Range("E19").Value = '  +0.09%  '
# Row 20
$ws.Range("D20").Value = '''7.41'
# Row 22
$ws.Range("E22").Value = '  +0.22%  '
# Row 23
$ws.Range("E23").Value = '  -4.95%  '
# Row 24
$ws.Range("D24").Value = '''9.15'
$ws.Range("E24").Value = '  -1.86%  '
# Row 25
$ws.Range("D25").Value = '''2.03'
$ws.Range("E25").Value = '  +4.59%  '
# Row 26
$ws.Range("D26").Value = '''151.47'
$ws.Range("E26").Value = '  -0.48%  '
# Row 27
$ws.Range("D27").Value = '''15.02'
$ws.Range("E27").Value = '  -1.36%  '
# Row 28
$ws.Range("E28").Value = '  -1.70%  '
# Row 29
$ws.Range("E29").Value = '  -1.84%  '
# Row 30
$ws.Range("E30").Value = '  +0.20%  '
# Row 31
$ws.Range("E31").Value = '  -1.84%  '
# Row 32
$ws.Range("D32").Value = '''0.0465'
$ws.Range("E32").Value = '  -1.91%  '
# Row 34
$ws.Range("E34").Value = '  +0.72%  '
# Row 35
$ws.Range("D35").Value = '1.397.35'
$ws.Range("E35").Value = '  -1.36%  '
# Row 36
$ws.Range("E36").Value = '  -2.45%  '
# Row 37
$ws.Range("E37").Value = '  -2.35%  '
# Row 38
$ws.Range("E38").Value = '  +1.76%  '
# Row 39
$ws.Range("D39").Value = '''2.60'
$ws.Range("E39").Value = '  +4.63%  '
# Row 40
$ws.Range("E40").Value = '  -0.45%  '
# Row 41
$ws.Range("E41").Value = '  -1.85%  '
# Row 42
$ws.Range("E42").Value = '  +0.22%  '
# Row 43
$ws.Range("D43").Value = '''0.796'
$ws.Range("E43").Value = '  -1.25%  '
# Row 44
$ws.Range("E44").Value = '  -1.02%  '
# Row 45
$ws.Range("D45").Value = '''1.87'
$ws.Range("E45").Value = '  +3.09%  '
# Row 46
$ws.Range("E46").Value = '  +0.71%  '
# Row 47
$ws.Range("E47").Value = '  -1.26%  '
# Row 48
$ws.Range("D48").Value = '1.717.30'
$ws.Range("E48").Value = '  +0.24%  '
# Row 49
$ws.Range("D49").Value = '''86.12'
$ws.Range("E49").Value = '  -0.99%  '
# Row 50
$ws.Range("E50").Value = '  -1.81%  '
# Row 51
$ws.Range("E51").Value = '  -1.31%  '
